$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (secrets / identifiers rotated for new Azure Key Vault names) ---
$ws.Range("B1").Value = "2YtQ=.R9kGf3yZk1xF.U=:=Fe[4:@vil"

$ws.Range("B2").Value = "https://igdcicd2.crm8.dynamics.com"
$ws.Range("C2").Value = "Source D365 Instance URL"

$ws.Range("B6").Value = "horrkvma7tdpvunkaaki5erfcf2hkkb4sl4tb2k37ciqdovsh7zq"

$ws.Range("B7").Value = "igdsa36@IGDCRM.onmicrosoft.com"

$ws.Range("B8").Value = "6d12e9fd-d509-4a1d-babf-40f344202c2b"

$ws.Range("B9").Value = "d9a1b506-a006-4359-966b-696cb2dad64d"

$ws.Range("B12").Value = "https://dd365key.vault.azure.net/"

# --- Hyperlinks: CRMSourceInstanceUrl (B2) and CRMSourceUserName (B4) become clickable ---
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:igdsa36@IGDCRM.onmicrosoft.com", "", "", "igdsa36@IGDCRM.onmicrosoft.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://igdcicd2.crm8.dynamics.com", "", "", "https://igdcicd2.crm8.dynamics.com")

# --- B6 (GitPassword secret) gets a dedicated font + centers vertically, no wrap ---
$ws.Range("B6").Font.Name = "Arial Unicode MS"
$ws.Range("B6").WrapText = $false
$ws.Range("B6").VerticalAlignment = -4108

# --- Row 2 height settles to the single-line height once the long URL got shortened ---
$ws.Rows(2).RowHeight = 15.75

# --- Selection moved to B7, and the view is no longer scrolled to row 7 ---
$ws.Range("B7").Select()

Write-Host "edit complete"
